# Auto-generated: apply scheduled market-price refresh to Mateus_Profits sheets.
# Updates cached currentAveragePrice / Leve-profit figures per sheet; a handful of
# rows also gain or lose a profit cell because the NQ/HQ price comparison flipped
# which side is cheaper (ClearContents removes the no-longer-applicable cell).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 269.3611
$ws.Range("I12").Value = 134.22858
$ws.Range("K12").Value = 134.22858
$ws.Range("M12").Value = 35.77142000000001
$ws.Range("H33").Value = 677.8889
$ws.Range("J33").Value = 879
$ws.Range("L33").Value = 879
$ws.Range("N33").Value = -1337
$ws.Range("H40").Value = 4180.275
$ws.Range("I40").Value = 2444.5833
$ws.Range("K40").Value = 2444.5833
$ws.Range("M40").Value = -2269.5833
$ws.Range("H55").Value = 260.05554
$ws.Range("I55").Value = 113.5
$ws.Range("J55").Value = 333.33334
$ws.Range("K55").Value = 113.5
$ws.Range("L55").Value = 333.33334
$ws.Range("M55").Value = 100.5
$ws.Range("N55").Value = -761.33334
$ws.Range("H58").Value = 265.8
$ws.Range("I58").Value = 184.33333
$ws.Range("K58").Value = 552.99999
$ws.Range("M58").Value = -402.99999
$ws.Range("H87").Value = 48031.332
$ws.Range("J87").Value = 48031.332
$ws.Range("L87").Value = 48031.332
$ws.Range("N87").Value = -50527.332
$ws.Range("H90").Value = 48031.332
$ws.Range("J90").Value = 48031.332
$ws.Range("L90").Value = 144093.996
$ws.Range("N90").Value = -156573.996
$ws.Range("H116").Value = 3374.875
$ws.Range("I116").Value = 3119.8
$ws.Range("J116").Value = 3800
$ws.Range("K116").Value = 3119.8
$ws.Range("L116").Value = 3800
$ws.Range("M116").Value = 322.1999999999998
$ws.Range("N116").Value = -10684
$ws.Range("H118").Value = 420
$ws.Range("I118").Value = 420
$ws.Range("K118").Value = 1260
$ws.Range("M118").Value = 397
$ws.Range("H137").Value = 1767.6
$ws.Range("J137").Value = 2582
$ws.Range("L137").Value = 7746
$ws.Range("N137").Value = -12846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 500.2857
$ws.Range("I5").Value = 130.75
$ws.Range("K5").Value = 130.75
$ws.Range("M5").Value = -18.75
$ws.Range("H26").Value = 1564.2
$ws.Range("I26").Value = 1564.2
$ws.Range("K26").Value = 1564.2
$ws.Range("M26").Value = -1234.2
$ws.Range("H32").Value = 5860.278
$ws.Range("I32").Value = 5860.278
$ws.Range("K32").Value = 5860.278
$ws.Range("M32").Value = -5573.278
$ws.Range("H74").Value = 3911.9285
$ws.Range("I74").Value = 2901.56
$ws.Range("K74").Value = 2901.56
$ws.Range("M74").Value = -2027.56
$ws.Range("H77").Value = 3911.9285
$ws.Range("I77").Value = 2901.56
$ws.Range("K77").Value = 14507.8
$ws.Range("M77").Value = -10139.8
$ws.Range("H80").Value = 31100
$ws.Range("H83").Value = 31100
$ws.Range("H115").Value = 40999
$ws.Range("I115").Value = 40998
$ws.Range("J115").Value = 41000
$ws.Range("K115").Value = 40998
$ws.Range("L115").Value = 41000
$ws.Range("M115").Value = -39431
$ws.Range("N115").Value = -44134
$ws.Range("H132").Value = 1427.48
$ws.Range("I132").Value = 1449.9
$ws.Range("J132").Value = 1337.8
$ws.Range("K132").Value = 4349.700000000001
$ws.Range("L132").Value = 4013.4
$ws.Range("M132").Value = -1819.700000000001
$ws.Range("N132").Value = -9073.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 500.2857
$ws.Range("I4").Value = 130.75
$ws.Range("K4").Value = 130.75
$ws.Range("M4").Value = -15.75
$ws.Range("H22").Value = 1599.6666
$ws.Range("I22").Value = 1619.6
$ws.Range("K22").Value = 1619.6
$ws.Range("M22").Value = -1446.6
$ws.Range("H35").Value = 32459.2
$ws.Range("J35").Value = 53432
$ws.Range("L35").Value = 53432
$ws.Range("N35").Value = -54052
$ws.Range("H82").Value = 21250
$ws.Range("J82").Value = 36500
$ws.Range("L82").Value = 36500
$ws.Range("N82").Value = -37266
$ws.Range("H85").Value = 21250
$ws.Range("J85").Value = 36500
$ws.Range("L85").Value = 36500
$ws.Range("N85").Value = -39152

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1039.9333
$ws.Range("I22").Value = 1181.7273
$ws.Range("K22").Value = 1181.7273
$ws.Range("M22").Value = -831.7273
$ws.Range("H31").Value = 31255454
$ws.Range("I31").Value = 71432056
$ws.Range("J31").Value = 6988.722
$ws.Range("K31").Value = 71432056
$ws.Range("L31").Value = 6988.722
$ws.Range("M31").Value = -71431761
$ws.Range("N31").Value = -7578.722
$ws.Range("H34").Value = 31255454
$ws.Range("I34").Value = 71432056
$ws.Range("J34").Value = 6988.722
$ws.Range("K34").Value = 71432056
$ws.Range("L34").Value = 6988.722
$ws.Range("M34").Value = -71431854
$ws.Range("N34").Value = -7392.722
$ws.Range("H41").Value = 21324.5
$ws.Range("J41").Value = 21324.5
$ws.Range("L41").Value = 21324.5
$ws.Range("N41").Value = -22180.5
$ws.Range("H51").Value = 22600
$ws.Range("J51").Value = 22600
$ws.Range("L51").Value = 22600
$ws.Range("N51").Value = -24072
$ws.Range("H61").Value = 22600
$ws.Range("J61").Value = 22600
$ws.Range("L61").Value = 22600
$ws.Range("N61").Value = -23296
$ws.Range("H74").Value = 41300
$ws.Range("J74").Value = 41300
$ws.Range("L74").Value = 41300
$ws.Range("N74").Value = -43048
$ws.Range("H77").Value = 41300
$ws.Range("J77").Value = 41300
$ws.Range("L77").Value = 123900
$ws.Range("N77").Value = -132636
$ws.Range("H107").Value = 302.66666
$ws.Range("I107").Value = 365.375
$ws.Range("J107").Value = 177.25
$ws.Range("K107").Value = 365.375
$ws.Range("L107").Value = 177.25
$ws.Range("M107").Value = 1554.625
$ws.Range("N107").Value = -4017.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1967.6111
$ws.Range("I2").Value = 27
$ws.Range("J2").Value = 8759.75
$ws.Range("K2").Value = 162
$ws.Range("L2").Value = 52558.5
$ws.Range("M2").Value = -49
$ws.Range("N2").Value = -52784.5
$ws.Range("H7").Value = 400297.8
$ws.Range("I7").Value = 500322.25
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 1500966.75
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -1500854.75
$ws.Range("N7").Value = -824
$ws.Range("H56").Value = 44571.285
$ws.Range("I56").Value = 44571.285
$ws.Range("K56").Value = 44571.285
$ws.Range("M56").Value = -44041.285
$ws.Range("H61").Value = 162.14285
$ws.Range("I61").Value = 44
$ws.Range("J61").Value = 250.75
$ws.Range("K61").Value = 132
$ws.Range("L61").Value = 752.25
$ws.Range("M61").Value = 83
$ws.Range("N61").Value = -1182.25
$ws.Range("H131").Value = 31252942
$ws.Range("J131").Value = 5036.2856
$ws.Range("L131").Value = 15108.8568
$ws.Range("N131").Value = -25188.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 9900
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("H132").Value = 2597.476
$ws.Range("I132").Value = 1836.1111
$ws.Range("J132").Value = 7165.6665
$ws.Range("K132").Value = 5508.3333
$ws.Range("L132").Value = 21496.9995
$ws.Range("M132").Value = -2978.3333
$ws.Range("N132").Value = -26556.9995
$ws.Range("M18").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1559.7
$ws.Range("J22").Value = 1766.6666
$ws.Range("L22").Value = 1766.6666
$ws.Range("N22").Value = -2356.6666
$ws.Range("H27").Value = 1559.7
$ws.Range("J27").Value = 1766.6666
$ws.Range("L27").Value = 1766.6666
$ws.Range("N27").Value = -1980.6666
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("H45").Value = 7249.5
$ws.Range("I45").Value = 7249.5
$ws.Range("K45").Value = 7249.5
$ws.Range("M45").Value = -6842.5
$ws.Range("H93").Value = 15374.667
$ws.Range("I93").Value = 2438.889
$ws.Range("J93").Value = 34778.332
$ws.Range("K93").Value = 2438.889
$ws.Range("L93").Value = 34778.332
$ws.Range("M93").Value = -1190.889
$ws.Range("N93").Value = -37274.332
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 15329.5
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("H15").Value = 7266
$ws.Range("J15").Value = 7266
$ws.Range("L15").Value = 7266
$ws.Range("N15").Value = -7842
$ws.Range("H54").Value = 23150
$ws.Range("J54").Value = 23150
$ws.Range("L54").Value = 23150
$ws.Range("N54").Value = -24190
$ws.Range("H107").Value = 1206.875
$ws.Range("I107").Value = 1155.8
$ws.Range("J107").Value = 1292
$ws.Range("K107").Value = 3467.4
$ws.Range("L107").Value = 3876
$ws.Range("M107").Value = -1547.4
$ws.Range("N107").Value = -7716
$ws.Range("H136").Value = 1693.65
$ws.Range("J136").Value = 5498.8
$ws.Range("L136").Value = 16496.4
$ws.Range("N136").Value = -21596.4
$ws.Range("N3").ClearContents()

Write-Output "Applied scheduled price refresh to Mateus_Profits workbook."
